$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper (scratch) cell used to write numeric-looking values as genuine text
# without Excel auto-converting them to numbers and without minting a new
# number-format style on the target cell (PasteSpecial xlPasteValues only
# carries the value + data type, not the source formatting/style).
$helper = $ws.Range("ZZ1")

$ws.Range('D2').Value = '25.965.78'
$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('D3').Value = '1.650.78'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  -0.14%  '
$helper.Value = '''216.72'
$helper.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -1.28%  '
$helper.Value = '''0.5184'
$helper.Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  -0.06%  '
$helper.Value = '''0.2618'
$helper.Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  -1.28%  '
$helper.Value = '''0.06239'
$helper.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  -1.40%  '
$helper.Value = '''20.50'
$helper.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -4.05%  '
$helper.Value = '''0.07720'
$helper.Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  -0.52%  '
$helper.Value = '''4.456'
$helper.Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.652.97'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').Value = '1.887.92'
$helper.Value = '''0.5411'
$helper.Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  -1.26%  '
$ws.Range('D16').Value = '0.0₅8062'
$ws.Range('E16').Value = '  -2.11%  '
$helper.Value = '''64.63'
$helper.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '26.010.88'
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('E19').Value = '  -0.19%  '
$helper.Value = '''4.561'
$helper.Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  -2.85%  '
$helper.Value = '''190.81'
$helper.Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  -0.44%  '
$helper.Value = '''9.974'
$helper.Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  -2.20%  '
$ws.Range('E23').Value = '  -0.17%  '
$helper.Value = '''5.965'
$helper.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  -4.04%  '
$helper.Value = '''137.93'
$helper.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  -0.77%  '
$helper.Value = '''0.1228'
$helper.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  -2.06%  '
$helper.Value = '''7.211'
$helper.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -1.31%  '
$helper.Value = '''16.05'
$helper.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  -0.28%  '
$helper.Value = '''1.401'
$helper.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -1.32%  '
$helper.Value = '''0.05900'
$helper.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  -2.58%  '
$ws.Range('E31').Value = '  -1.21%  '
$helper.Value = '''3.513'
$helper.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  -1.21%  '
$helper.Value = '''3.241'
$helper.Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  -3.95%  '
$ws.Range('E34').Value = '  -6.08%  '
$helper.Value = '''2.417'
$helper.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  -0.44%  '
$helper.Value = '''0.9456'
$helper.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  -4.36%  '
$helper.Value = '''2.752'
$helper.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -0.74%  '
$helper.Value = '''0.5629'
$helper.Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  -5.90%  '
$ws.Range('E39').Value = '  -0.54%  '
$helper.Value = '''5.852'
$helper.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  -2.20%  '
$helper.Value = '''0.8458'
$helper.Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  -0.66%  '
$helper.Value = '''1.004'
$helper.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -0.07%  '
$helper.Value = '''100.75'
$helper.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  +0.75%  '
$helper.Value = '''998.09'
$helper.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  -6.80%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈107'
$ws.Range('E46').Value = '  -2.44%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$helper.Value = '''56.32'
$helper.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  -2.14%  '
$helper.Value = '''0.9986'
$helper.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$helper.Value = '''7.974'
$helper.Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$helper.Value = '''0.4309'
$helper.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  +1.81%  '
$helper.Value = '''0.05151'
$helper.Copy()
$ws.Range('D51').PasteSpecial(-4163)

$helper.Clear()
$excel.CutCopyMode = $false

